$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "76.370.35"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.51%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.045.75"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +3.08%  "

$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "198.80"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.46%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "618.54"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.97%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.550"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.31%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.210"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +7.15%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "3.044.79"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.01%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.442"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.95%  "

$ws.Range("E12").Value = "  -0.17%  "

$ws.Range("E13").Value = "  +6.51%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.603.92"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.71%  "

$ws.Range("E15").Value = "  +2.58%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "76.239.66"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.69%  "

$ws.Range("E17").Value = "  +2.81%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.048.05"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +4.09%  "

$ws.Range("E19").Value = "  +1.74%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.91"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.22%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "380.94"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.59%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.45"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +8.63%  "

$ws.Range("E23").Value = "  +2.25%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.200.82"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.18%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "72.41"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.08%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.03%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "4.33"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.95%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.83"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.76%  "

$ws.Range("E29").Value = "  +1.62%  "

$ws.Range("E30").Value = "  -0.22%  "

$ws.Range("E31").Value = "  +1.36%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.39"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.48%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "495.43"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.02%  "

$ws.Range("E34").Value = "  +5.03%  "

$ws.Range("E35").Value = "  -0.11%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.124"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +12.77%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "20.67"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.66%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "162.32"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.48%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "20.06"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.42%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "191.72"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +5.74%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.377"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.71%  "

$ws.Range("E42").Value = "  -8.32%  "

$ws.Range("E43").Value = "  -0.04%  "

$ws.Range("E44").Value = "  +4.30%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.783"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +18.48%  "

$ws.Range("E46").Value = "  +5.82%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "41.21"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.75%  "

$ws.Range("E48").Value = "  +0.13%  "

$ws.Range("E49").Value = "  +5.90%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.594"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.52%  "

$ws.Range("E51").Value = "  -0.15%  "
